$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was reported; insert a row before row 435
# (shifting the existing Femacal de La Calera - Ajo rows down by one)
# and populate it with the new observation.
$ws.Rows.Item(435).Insert()

$ws.Cells.Item(435, 1).Value = 3
$ws.Cells.Item(435, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(435, 3).Value = "Coquimbo"
$ws.Cells.Item(435, 4).Value = 44746
$ws.Cells.Item(435, 5).Value = 5
$ws.Cells.Item(435, 6).Value = 100112003
$ws.Cells.Item(435, 7).Value = "Ajo"
$ws.Cells.Item(435, 8).Value = "Chino"
$ws.Cells.Item(435, 9).Value = "Primera"
$ws.Cells.Item(435, 10).Value = 85
$ws.Cells.Item(435, 11).Value = 17000
$ws.Cells.Item(435, 12).Value = 17500
$ws.Cells.Item(435, 13).Value = 17265
$ws.Cells.Item(435, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(435, 15).Value = "China"
$ws.Cells.Item(435, 16).Value = 1726
$ws.Cells.Item(435, 17).Value = 10
$ws.Cells.Item(435, 18).Value = "Hortaliza"
